# Auto-generated script applying the Sargatanas_Profits.xlsx diff
# (scheduled market-data refresh: updates currentAveragePrice* / Leve* columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4548.4
$ws.Range("I76").Value = 4330.3335
$ws.Range("J76").Value = 4875.5
$ws.Range("K76").Value = 4330.3335
$ws.Range("L76").Value = 4875.5
$ws.Range("M76").Value = -4015.3335
$ws.Range("N76").Value = -5505.5
$ws.Range("H79").Value = 4548.4
$ws.Range("I79").Value = 4330.3335
$ws.Range("J79").Value = 4875.5
$ws.Range("K79").Value = 4330.3335
$ws.Range("L79").Value = 4875.5
$ws.Range("M79").Value = -3238.3335
$ws.Range("N79").Value = -7059.5
$ws.Range("H112").Value = 6641.2812
$ws.Range("J112").Value = 6641.2812
$ws.Range("L112").Value = 19923.8436
$ws.Range("N112").Value = -22139.8436
$ws.Range("H116").Value = 9800.200000000001
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 9800.200000000001
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 9800.200000000001
$ws.Range("M116").Value = ""
$ws.Range("N116").Value = -16684.2
$ws.Range("H132").Value = 1458.5416
$ws.Range("I132").Value = 1213.3243
$ws.Range("J132").Value = 2283.3635
$ws.Range("K132").Value = 3639.9729
$ws.Range("L132").Value = 6850.0905
$ws.Range("M132").Value = -1109.9729
$ws.Range("N132").Value = -11910.0905
$ws.Range("H138").Value = 3715.7021
$ws.Range("I138").Value = 923.45
$ws.Range("J138").Value = 5784.037
$ws.Range("K138").Value = 2770.35
$ws.Range("L138").Value = 17352.111
$ws.Range("M138").Value = 2369.65
$ws.Range("N138").Value = -27632.111
$ws.Range("H141").Value = 2101.8333
$ws.Range("I141").Value = 2108
$ws.Range("K141").Value = 6324
$ws.Range("M141").Value = -1144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3995.3865
$ws.Range("I61").Value = 2327.8333
$ws.Range("K61").Value = 2327.8333
$ws.Range("M61").Value = -2115.8333
$ws.Range("H97").Value = 8565.6
$ws.Range("I97").Value = 615.63635
$ws.Range("J97").Value = 30428
$ws.Range("K97").Value = 615.63635
$ws.Range("L97").Value = 30428
$ws.Range("M97").Value = -119.63635
$ws.Range("N97").Value = -31420
$ws.Range("H132").Value = 5473.909
$ws.Range("I132").Value = 4087.5898
$ws.Range("K132").Value = 12262.7694
$ws.Range("M132").Value = -9732.769400000001
$ws.Range("H136").Value = 3995.3865
$ws.Range("I136").Value = 2327.8333
$ws.Range("K136").Value = 6983.499899999999
$ws.Range("M136").Value = -4433.499899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5324.773
$ws.Range("I132").Value = 1976.6666
$ws.Range("J132").Value = 12499.286
$ws.Range("K132").Value = 5929.9998
$ws.Range("L132").Value = 37497.858
$ws.Range("M132").Value = -3399.9998
$ws.Range("N132").Value = -42557.858
$ws.Range("H137").Value = 45000
$ws.Range("J137").Value = 45000
$ws.Range("L137").Value = 45000
$ws.Range("N137").Value = -55200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 2000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 2000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 6000
$ws.Range("M44").Value = ""
$ws.Range("N44").Value = -6796
$ws.Range("H80").Value = 40004240
$ws.Range("I80").Value = 26319736
$ws.Range("J80").Value = 83338500
$ws.Range("K80").Value = 78959208
$ws.Range("L80").Value = 250015500
$ws.Range("M80").Value = -78958272
$ws.Range("N80").Value = -250017372
$ws.Range("H83").Value = 40004240
$ws.Range("I83").Value = 26319736
$ws.Range("J83").Value = 83338500
$ws.Range("K83").Value = 236877624
$ws.Range("L83").Value = 750046500
$ws.Range("M83").Value = -236872944
$ws.Range("N83").Value = -750055860
$ws.Range("H92").Value = 6994823
$ws.Range("J92").Value = 6994823
$ws.Range("L92").Value = 20984469
$ws.Range("N92").Value = -20986965
$ws.Range("H122").Value = 3146252
$ws.Range("J122").Value = 8333
$ws.Range("L122").Value = 74997
$ws.Range("N122").Value = -79897
$ws.Range("H132").Value = 13902.471
$ws.Range("I132").Value = 5557.143
$ws.Range("K132").Value = 50014.287
$ws.Range("M132").Value = -47484.287
$ws.Range("H138").Value = 80324.38
$ws.Range("I138").Value = 86601.414
$ws.Range("J138").Value = 5000
$ws.Range("K138").Value = 259804.242
$ws.Range("L138").Value = 15000
$ws.Range("M138").Value = -254664.242
$ws.Range("N138").Value = -25280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2095.1516
$ws.Range("J97").Value = 3246
$ws.Range("L97").Value = 3246
$ws.Range("N97").Value = -4238
$ws.Range("H102").Value = 2957.3142
$ws.Range("I102").Value = 2541.1482
$ws.Range("J102").Value = 4361.875
$ws.Range("K102").Value = 2541.1482
$ws.Range("L102").Value = 4361.875
$ws.Range("M102").Value = -919.1482000000001
$ws.Range("N102").Value = -7605.875
$ws.Range("H122").Value = 78638.92999999999
$ws.Range("I122").Value = 130257.25
$ws.Range("K122").Value = 390771.75
$ws.Range("M122").Value = -388321.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6190.3076
$ws.Range("I7").Value = 4391.1665
$ws.Range("J7").Value = 7732.4287
$ws.Range("K7").Value = 4391.1665
$ws.Range("L7").Value = 7732.4287
$ws.Range("M7").Value = -4279.1665
$ws.Range("N7").Value = -7956.4287
$ws.Range("H22").Value = 5332
$ws.Range("J22").Value = 6109.3335
$ws.Range("L22").Value = 6109.3335
$ws.Range("N22").Value = -6699.3335
$ws.Range("H27").Value = 5332
$ws.Range("J27").Value = 6109.3335
$ws.Range("L27").Value = 6109.3335
$ws.Range("N27").Value = -6323.3335
$ws.Range("I55").Value = 250000240
$ws.Range("J55").Value = 791.75
$ws.Range("K55").Value = 250000240
$ws.Range("L55").Value = 791.75
$ws.Range("M55").Value = -250000067
$ws.Range("N55").Value = -1137.75
$ws.Range("H126").Value = 6190.3076
$ws.Range("I126").Value = 4391.1665
$ws.Range("J126").Value = 7732.4287
$ws.Range("K126").Value = 13173.4995
$ws.Range("L126").Value = 23197.2861
$ws.Range("M126").Value = -10703.4995
$ws.Range("N126").Value = -28137.2861
$ws.Range("H132").Value = 17866210
$ws.Range("I132").Value = 71431704
$ws.Range("J132").Value = 11047.333
$ws.Range("K132").Value = 214295112
$ws.Range("L132").Value = 33141.999
$ws.Range("M132").Value = -214292582
$ws.Range("N132").Value = -38201.999
$ws.Range("H136").Value = 9204.773999999999
$ws.Range("I136").Value = 2891
$ws.Range("J136").Value = 11787.682
$ws.Range("K136").Value = 8673
$ws.Range("L136").Value = 35363.046
$ws.Range("M136").Value = -6123
$ws.Range("N136").Value = -40463.046

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 17666.666
$ws.Range("J43").Value = 18800
$ws.Range("L43").Value = 18800
$ws.Range("N43").Value = -19098
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = ""
$ws.Range("H54").Value = 14076.923
$ws.Range("H61").Value = 10051
$ws.Range("I61").Value = 10051
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 10051
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -9759
$ws.Range("N61").Value = ""
$ws.Range("H126").Value = 4067.7778
$ws.Range("I126").Value = 2944.2856
$ws.Range("K126").Value = 8832.856800000001
$ws.Range("M126").Value = -6362.856800000001
$ws.Range("H132").Value = 19246978
$ws.Range("I132").Value = 27785528
$ws.Range("J132").Value = 35236.25
$ws.Range("K132").Value = 83356584
$ws.Range("L132").Value = 105708.75
$ws.Range("M132").Value = -83354054
$ws.Range("N132").Value = -110768.75
$ws.Range("H136").Value = 52687370
$ws.Range("I136").Value = 100000670
$ws.Range("K136").Value = 300002010
$ws.Range("M136").Value = -299999460

